# Apply the "added checked invoices persists" edit:
# Swap the two invoice date strings held in D1/D2 on the first sheet.
# Before: D1 -> "30.08.2016", D2 -> "30.09.2016"
# After:  D1 -> "30.10.2016", D2 -> "30.09.2016"
# (D1 keeps its quote-prefixed text style, so re-enter it with a leading
# apostrophe to force Excel to treat it as literal text, same as before.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D1").Value = "'30.10.2016"
$ws.Range("D2").Value = "30.09.2016"

# Move/restore the sheet's active selection from E3 to D2.
$ws.Range("D2").Select()
